$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Mark the "Milestone Completed" columns (E = milestone numeral, F = confirmation "X")
# for the rows that were completed in this commit.
$rows = @(7, 21, 39, 40, 57)
foreach ($r in $rows) {
    $ws.Range("F$r").Value = "X"
    $ws.Range("E$r").Value = "I"
}

# Update the active selection left on the sheet to E41, matching the author's cursor position.
$ws.Range("E41").Select()

$wb.Save()
